$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New values for cells X6 and Y6 in the existing row 6
$ws.Range("X6").Value = 0.059999000000001246
$ws.Range("Y6").Value = "Up"

# New row 7 data
$ws.Range("A7").Value = 42648.886597222219
$ws.Range("B7").Value = 11
$ws.Range("C7").Value = "Buy"
$ws.Range("D7").Value = 38
$ws.Range("E7").Value = 15151
$ws.Range("F7").Value = 1666
$ws.Range("G7").Value = 65
$ws.Range("H7").Value = 31
$ws.Range("I7").Value = 91
$ws.Range("J7").Value = 8
$ws.Range("K7").Value = 13409
$ws.Range("L7").Value = 298
$ws.Range("M7").Value = 142
$ws.Range("N7").Value = 89
$ws.Range("O7").Value = 8
$ws.Range("P7").Value = "Bag"
$ws.Range("Q7").Value = 41.162214763508182
$ws.Range("R7").Value = 0
$ws.Range("S7").Value = 0.061600000000000002
$ws.Range("T7").Value = -0.032000000000000001
$ws.Range("U7").Value = 2.2599999999999998
$ws.Range("V7").Value = "N/A"
$ws.Range("W7").Value = 0

# Apply number formats matching column A (date) and S/T (percentage) styles
$ws.Range("A7").NumberFormat = "m/d/yy h:mm"
$ws.Range("S7").NumberFormat = $ws.Range("S6").NumberFormat
$ws.Range("T7").NumberFormat = $ws.Range("T6").NumberFormat
